# Fruta / hortaliza, semanal
# A new weekly price record is inserted as row 17 (shifting the existing
# rows 17-79 down to 18-80). The new row reuses the same market/product
# metadata as the row that follows it, with its own date/volume/price
# figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 17; rows 17:79 shift down to 18:80.
$ws.Rows("17:17").Insert()

# Seed the new row 17 from (the now-shifted) row 18, which carries the
# shared Mercado/Producto/Variedad/Calidad/Unidad/Origen metadata.
$ws.Range("A18:T18").Copy()
$ws.Range("A17:T17").PasteSpecial()

# Overwrite the new row's own date, volume and price figures.
$ws.Range("D17").Value2 = 44690
$ws.Range("M17").Value2 = 25
$ws.Range("N17").Value2 = 25000
$ws.Range("O17").Value2 = 25000
$ws.Range("P17").Value2 = 25000
$ws.Range("S17").Value2 = 2500
